$wb = $excel.ActiveWorkbook

# --- Overview sheet: Latest Handoff Date column (D) ---
$ws = $wb.Worksheets.Item("Overview")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $cell = $ws.Range("D$r")
    $cell.NumberFormat = "General"
    $cell.Value = "2016-19-11 12:19:41"
}

# --- zh-cn sheet: Latest Handoff Datetime column (E) ---
$ws = $wb.Worksheets.Item("zh-cn")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $cell = $ws.Range("E$r")
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $cell.Value = "2016-03-11 12:19:37"
}

# --- de-de sheet: Latest Handoff Datetime column (E) ---
$ws = $wb.Worksheets.Item("de-de")
foreach ($r in 7,10,11,12,13,14,15,16) {
    $cell = $ws.Range("E$r")
    $cell.NumberFormat = "yyyy-mm-dd HH:mm:ss"
    $cell.Value = "2016-03-11 12:19:41"
}
